$d = $word.ActiveDocument

$replacements = @(
    @("520÷5=", "623÷8="),
    @("785÷3=", "756÷5="),
    @("197÷9=", "730÷8="),
    @("957÷3=", "397÷8="),
    @("956÷8=", "938÷3="),
    @("625÷5=", "377÷3="),
    @("760÷9=", "964÷8="),
    @("711÷8=", "516÷9="),
    @("139÷7=", "259÷8="),
    @("224÷7=", "711÷5="),
    @("669÷9=", "544÷5="),
    @("598÷4=", "747÷2="),
    @("792÷8=", "534÷8="),
    @("741÷3=", "895÷6="),
    @("228÷8=", "900÷3="),
    @("296÷3=", "920÷5="),
    @("450÷4=", "449÷5="),
    @("989÷6=", "854÷3="),
    @("516÷3=", "463÷9="),
    @("939÷3=", "625÷8="),
    @("462÷4=", "350÷7="),
    @("703÷4=", "156÷3="),
    @("118÷3=", "192÷2="),
    @("147÷8=", "685÷9="),
    @("106÷2=", "360÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
